# Rename the worksheet tab from "Mapping Tag Glossary (2)" to
# "Mapping Tag Glossary" (commit: "Update tab name in this BRIDG Mapping
# Tag Glossary spreadsheet").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Mapping Tag Glossary"
